$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Save a snapshot of the current (pre-edit) values for the cells that
# participate in the cyclic shift: row 2 -> row 3 -> row 4 -> row 5 -> row 2
$cols = @("D","J","K","L","M","P")

$row2 = @{}
$row3 = @{}
$row4 = @{}
$row5 = @{}

foreach ($col in $cols) {
    $row2[$col] = $ws.Range("${col}2").Value2
    $row3[$col] = $ws.Range("${col}3").Value2
    $row4[$col] = $ws.Range("${col}4").Value2
    $row5[$col] = $ws.Range("${col}5").Value2
}

# Apply the cyclic shift:
#   new row2 <- old row5
#   new row3 <- old row2
#   new row4 <- old row3
#   new row5 <- old row4
foreach ($col in $cols) {
    $ws.Range("${col}2").Value2 = $row5[$col]
    $ws.Range("${col}3").Value2 = $row2[$col]
    $ws.Range("${col}4").Value2 = $row3[$col]
    $ws.Range("${col}5").Value2 = $row4[$col]
}
